$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading zeros must be preserved as text, so prefix values that look like
# numbers with an apostrophe to force text entry (keeps existing "General"
# number format / style instead of reformatting as a numeric value).
$ws.Range("B1").Value = "'00059171"
$ws.Range("F1").Value = "Christopher Felski"

$ws.Range("B3").Value = "'300006839"
$ws.Range("F3").Value = "Christopher Felski"

$ws.Range("B5").Value = "'0880011949"
$ws.Range("F5").Value = "01/26"

$ws.Range("E10").Value = "prop damage due to limb strike. props incorrectly replaced resulting in rollover crashes. this may have resulted in broken motor arm and gimble damage."
